$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is updated in place: it now carries the "new" exposure-time entry
# for the existing Point Cook venue (previously at row 3 with a different
# exposure period / "old" marker).
$ws.Range("B3").Value = "The Coffeeologist Cafe  70/300 Point Cook Rd  Point Cook VIC 3030"
$ws.Range("C3").Value = "11:30am - 12:10pm  0/2/2021"
$ws.Range("D3").Value = "Case attended venue"
$ws.Range("E3").Value = "new"

# Row 4: the original Point Cook venue text/time moves down, marked "old"
$ws.Range("A4").Value = "Point Cook"
$ws.Range("B4").Value = "The Coffeeologist Cafe, 70/300 Point Cook Rd , Point Cook VIC 3030"
$ws.Range("C4").Value = "11:00am - 11:40am 8/2/2021"
$ws.Range("D4").Value = "Case attended venue"
$ws.Range("E4").Value = "old"

# Row 5: the previous row-3 venue text/time, marked "old"
$ws.Range("A5").Value = "Point Cook"
$ws.Range("B5").Value = "The Coffeeologist Cafe, 70/300 Point Cook Rd, Point Cook VIC 3030"
$ws.Range("C5").Value = "11:30am - 12:10pm 10/2/2021"
$ws.Range("D5").Value = "Case attended venue"
$ws.Range("E5").Value = "old"

# Row 6: new Sunbury location, "new"
$ws.Range("A6").Value = "Sunbury"
$ws.Range("B6").Value = "Sunbury Square Shopping Centre  2-28 Evans St  Sunbury VIC 3429"
$ws.Range("C6").Value = "3:40pm - 4:30pm 5/2/2021"
$ws.Range("D6").Value = "Case attended venue"
$ws.Range("E6").Value = "new"

# Row 7: same Sunbury location with slightly different site text, "old"
$ws.Range("A7").Value = "Sunbury"
$ws.Range("B7").Value = "Sunbury Square Shopping Centre, 2-28 Evans street, Sunbury"
$ws.Range("C7").Value = "3:40pm - 4:30pm 5/2/2021"
$ws.Range("D7").Value = "Case attended venue"
$ws.Range("E7").Value = "old"

# Match the saved selection state recorded in the diff
$ws.Range("C3").Select()
